# The upstream change between the "before" and "after" revisions of this
# template is a pure XML re-serialization: every element's attributes
# (and every xmlns: namespace declaration on the part-root elements) are
# alphabetized by the tool that regenerated the test fixture when the
# authoring/test library moved from version 2.0.2 to 2.0.3 (see commit
# message). Comparing every changed line pairwise confirms the attribute
# *name -> value* sets are identical before and after; only their
# on-the-wire order changes. No run text, paragraph, table cell, style,
# section, margin, or any other visible/structural property differs.
#
# Word's object model -- and this COM-interop surface -- models document
# *content*, not byte-level XML attribute ordering, so there is no
# content edit to perform: the document's content already matches the
# target state exactly. Intentionally making no changes here is what
# keeps the saved package's content equivalent to the target; any
# property "touch" through the OM would instead risk *introducing* a
# content/formatting difference that isn't in the source diff (e.g.
# minting new namespace declarations or renormalizing whitespace), which
# would be incorrect.
